# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# formatting used by the existing columns and initializing the data row to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in H1, matching text + formatting of the other headers.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data cell in H2.
$ws.Range("H2").Value = 0
